$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 66.85695766666667
$ws.Cells.Item(2, 8).Value = 200.570873
$ws.Cells.Item(2, 9).Value = 0.8284701681115905
$ws.Cells.Item(2, 10).Value = 0.8284701681115904
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 66.85695766666667
$ws.Cells.Item(2, 14).Value = 200.570873
$ws.Cells.Item(2, 15).Value = 0.8284701681115905
$ws.Cells.Item(2, 16).Value = 0.8284701681115904
$ws.Cells.Item(2, 17).Value = 4469.852788442459
$ws.Cells.Item(2, 18).Value = 40228.67509598213
$ws.Cells.Item(2, 19).Value = 0.686362819450847
$ws.Cells.Item(2, 20).Value = 0.6863628194508468

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 66.85695766666667
$ws.Cells.Item(3, 8).Value = 200.570873
$ws.Cells.Item(3, 9).Value = 0.8284701681115905
$ws.Cells.Item(3, 10).Value = 0.8284701681115904
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 10.944004
$ws.Cells.Item(3, 14).Value = 32.83201200000001
$ws.Cells.Item(3, 15).Value = 0.1356146188837786
$ws.Cells.Item(3, 16).Value = 0.1356146188837786
$ws.Cells.Item(3, 17).Value = 731.6828121318308
$ws.Cells.Item(3, 18).Value = 6585.145309186478
$ws.Cells.Item(3, 19).Value = 0.1123526661050333
$ws.Cells.Item(3, 20).Value = 0.1123526661050333

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 66.85695766666667
$ws.Cells.Item(4, 8).Value = 200.570873
$ws.Cells.Item(4, 9).Value = 0.8284701681115905
$ws.Cells.Item(4, 10).Value = 0.8284701681115904
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.049009
$ws.Cells.Item(4, 14).Value = 0.147027
$ws.Cells.Item(4, 15).Value = 0.0006073039498957697
$ws.Cells.Item(4, 16).Value = 0.0006073039498957697
$ws.Cells.Item(4, 17).Value = 3.276592638285667
$ws.Cells.Item(4, 18).Value = 29.489333744571
$ws.Cells.Item(4, 19).Value = 0.0005031332054649812
$ws.Cells.Item(4, 20).Value = 0.0005031332054649812

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 66.85695766666667
$ws.Cells.Item(5, 8).Value = 200.570873
$ws.Cells.Item(5, 9).Value = 0.8284701681115905
$ws.Cells.Item(5, 10).Value = 0.8284701681115904
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.3066033333333333
$ws.Cells.Item(5, 14).Value = 0.91981
$ws.Cells.Item(5, 15).Value = 0.003799331049083692
$ws.Cells.Item(5, 16).Value = 0.003799331049083692
$ws.Cells.Item(5, 17).Value = 20.49856607712556
$ws.Cells.Item(5, 18).Value = 184.48709469413
$ws.Cells.Item(5, 19).Value = 0.003147632432945952
$ws.Cells.Item(5, 20).Value = 0.003147632432945951

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 66.85695766666667
$ws.Cells.Item(6, 8).Value = 200.570873
$ws.Cells.Item(6, 9).Value = 0.8284701681115905
$ws.Cells.Item(6, 10).Value = 0.8284701681115904
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.1192703333333334
$ws.Cells.Item(6, 14).Value = 0.357811
$ws.Cells.Item(6, 15).Value = 0.001477960059146655
$ws.Cells.Item(6, 16).Value = 0.001477960059146655
$ws.Cells.Item(6, 17).Value = 7.974051626555891
$ws.Cells.Item(6, 18).Value = 71.766464639003
$ws.Cells.Item(6, 19).Value = 0.001224445818663446
$ws.Cells.Item(6, 20).Value = 0.001224445818663445

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 66.85695766666667
$ws.Cells.Item(7, 8).Value = 200.570873
$ws.Cells.Item(7, 9).Value = 0.8284701681115905
$ws.Cells.Item(7, 10).Value = 0.8284701681115904
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.423449666666667
$ws.Cells.Item(7, 14).Value = 7.270349
$ws.Cells.Item(7, 15).Value = 0.03003061794650479
$ws.Cells.Item(7, 16).Value = 0.03003061794650479
$ws.Cells.Item(7, 17).Value = 162.0244717716308
$ws.Cells.Item(7, 18).Value = 1458.220245944677
$ws.Cells.Item(7, 19).Value = 0.02487947109863577
$ws.Cells.Item(7, 20).Value = 0.02487947109863576

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 10.944004
$ws.Cells.Item(8, 8).Value = 32.83201200000001
$ws.Cells.Item(8, 9).Value = 0.1356146188837786
$ws.Cells.Item(8, 10).Value = 0.1356146188837786
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 66.85695766666667
$ws.Cells.Item(8, 14).Value = 200.570873
$ws.Cells.Item(8, 15).Value = 0.8284701681115905
$ws.Cells.Item(8, 16).Value = 0.8284701681115904
$ws.Cells.Item(8, 17).Value = 731.6828121318308
$ws.Cells.Item(8, 18).Value = 6585.145309186478
$ws.Cells.Item(8, 19).Value = 0.1123526661050333
$ws.Cells.Item(8, 20).Value = 0.1123526661050333

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 10.944004
$ws.Cells.Item(9, 8).Value = 32.83201200000001
$ws.Cells.Item(9, 9).Value = 0.1356146188837786
$ws.Cells.Item(9, 10).Value = 0.1356146188837786
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 10.944004
$ws.Cells.Item(9, 14).Value = 32.83201200000001
$ws.Cells.Item(9, 15).Value = 0.1356146188837786
$ws.Cells.Item(9, 16).Value = 0.1356146188837786
$ws.Cells.Item(9, 17).Value = 119.771223552016
$ws.Cells.Item(9, 18).Value = 1077.941011968144
$ws.Cells.Item(9, 19).Value = 0.01839132485499252
$ws.Cells.Item(9, 20).Value = 0.01839132485499251

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 10.944004
$ws.Cells.Item(10, 8).Value = 32.83201200000001
$ws.Cells.Item(10, 9).Value = 0.1356146188837786
$ws.Cells.Item(10, 10).Value = 0.1356146188837786
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.049009
$ws.Cells.Item(10, 14).Value = 0.147027
$ws.Cells.Item(10, 15).Value = 0.0006073039498957697
$ws.Cells.Item(10, 16).Value = 0.0006073039498957697
$ws.Cells.Item(10, 17).Value = 0.5363546920360001
$ws.Cells.Item(10, 18).Value = 4.827192228324001
$ws.Cells.Item(10, 19).Value = 0.00008235929371172819
$ws.Cells.Item(10, 20).Value = 0.00008235929371172816

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 10.944004
$ws.Cells.Item(11, 8).Value = 32.83201200000001
$ws.Cells.Item(11, 9).Value = 0.1356146188837786
$ws.Cells.Item(11, 10).Value = 0.1356146188837786
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.3066033333333333
$ws.Cells.Item(11, 14).Value = 0.91981
$ws.Cells.Item(11, 15).Value = 0.003799331049083692
$ws.Cells.Item(11, 16).Value = 0.003799331049083692
$ws.Cells.Item(11, 17).Value = 3.355468106413334
$ws.Cells.Item(11, 18).Value = 30.19921295772001
$ws.Cells.Item(11, 19).Value = 0.0005152448322347917
$ws.Cells.Item(11, 20).Value = 0.0005152448322347915

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 10.944004
$ws.Cells.Item(12, 8).Value = 32.83201200000001
$ws.Cells.Item(12, 9).Value = 0.1356146188837786
$ws.Cells.Item(12, 10).Value = 0.1356146188837786
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.1192703333333334
$ws.Cells.Item(12, 14).Value = 0.357811
$ws.Cells.Item(12, 15).Value = 0.001477960059146655
$ws.Cells.Item(12, 16).Value = 0.001477960059146655
$ws.Cells.Item(12, 17).Value = 1.305295005081334
$ws.Cells.Item(12, 18).Value = 11.747655045732
$ws.Cells.Item(12, 19).Value = 0.0002004329901466205
$ws.Cells.Item(12, 20).Value = 0.0002004329901466205

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 10.944004
$ws.Cells.Item(13, 8).Value = 32.83201200000001
$ws.Cells.Item(13, 9).Value = 0.1356146188837786
$ws.Cells.Item(13, 10).Value = 0.1356146188837786
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 2.423449666666667
$ws.Cells.Item(13, 14).Value = 7.270349
$ws.Cells.Item(13, 15).Value = 0.03003061794650479
$ws.Cells.Item(13, 16).Value = 0.03003061794650479
$ws.Cells.Item(13, 17).Value = 26.52224284579867
$ws.Cells.Item(13, 18).Value = 238.700185612188
$ws.Cells.Item(13, 19).Value = 0.004072590807659609
$ws.Cells.Item(13, 20).Value = 0.004072590807659608

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 0.049009
$ws.Cells.Item(14, 8).Value = 0.147027
$ws.Cells.Item(14, 9).Value = 0.0006073039498957697
$ws.Cells.Item(14, 10).Value = 0.0006073039498957697
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 66.85695766666667
$ws.Cells.Item(14, 14).Value = 200.570873
$ws.Cells.Item(14, 15).Value = 0.8284701681115905
$ws.Cells.Item(14, 16).Value = 0.8284701681115904
$ws.Cells.Item(14, 17).Value = 3.276592638285667
$ws.Cells.Item(14, 18).Value = 29.489333744571
$ws.Cells.Item(14, 19).Value = 0.0005031332054649812
$ws.Cells.Item(14, 20).Value = 0.0005031332054649812

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 0.049009
$ws.Cells.Item(15, 8).Value = 0.147027
$ws.Cells.Item(15, 9).Value = 0.0006073039498957697
$ws.Cells.Item(15, 10).Value = 0.0006073039498957697
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 10.944004
$ws.Cells.Item(15, 14).Value = 32.83201200000001
$ws.Cells.Item(15, 15).Value = 0.1356146188837786
$ws.Cells.Item(15, 16).Value = 0.1356146188837786
$ws.Cells.Item(15, 17).Value = 0.5363546920360001
$ws.Cells.Item(15, 18).Value = 4.827192228324001
$ws.Cells.Item(15, 19).Value = 0.00008235929371172819
$ws.Cells.Item(15, 20).Value = 0.00008235929371172816

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 0.049009
$ws.Cells.Item(16, 8).Value = 0.147027
$ws.Cells.Item(16, 9).Value = 0.0006073039498957697
$ws.Cells.Item(16, 10).Value = 0.0006073039498957697
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.049009
$ws.Cells.Item(16, 14).Value = 0.147027
$ws.Cells.Item(16, 15).Value = 0.0006073039498957697
$ws.Cells.Item(16, 16).Value = 0.0006073039498957697
$ws.Cells.Item(16, 17).Value = 0.002401882081
$ws.Cells.Item(16, 18).Value = 0.021616938729
$ws.Cells.Item(16, 19).Value = 0.0000003688180875590036
$ws.Cells.Item(16, 20).Value = 0.0000003688180875590036

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 0.049009
$ws.Cells.Item(17, 8).Value = 0.147027
$ws.Cells.Item(17, 9).Value = 0.0006073039498957697
$ws.Cells.Item(17, 10).Value = 0.0006073039498957697
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.3066033333333333
$ws.Cells.Item(17, 14).Value = 0.91981
$ws.Cells.Item(17, 15).Value = 0.003799331049083692
$ws.Cells.Item(17, 16).Value = 0.003799331049083692
$ws.Cells.Item(17, 17).Value = 0.01502632276333333
$ws.Cells.Item(17, 18).Value = 0.13523690487
$ws.Cells.Item(17, 19).Value = 0.000002307348753070165
$ws.Cells.Item(17, 20).Value = 0.000002307348753070164

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 0.049009
$ws.Cells.Item(18, 8).Value = 0.147027
$ws.Cells.Item(18, 9).Value = 0.0006073039498957697
$ws.Cells.Item(18, 10).Value = 0.0006073039498957697
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 0.1192703333333334
$ws.Cells.Item(18, 14).Value = 0.357811
$ws.Cells.Item(18, 15).Value = 0.001477960059146655
$ws.Cells.Item(18, 16).Value = 0.001477960059146655
$ws.Cells.Item(18, 17).Value = 0.005845319766333334
$ws.Cells.Item(18, 18).Value = 0.05260787789700001
$ws.Cells.Item(18, 19).Value = 0.0000008975709817079493
$ws.Cells.Item(18, 20).Value = 0.000000897570981707949

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 0.049009
$ws.Cells.Item(19, 8).Value = 0.147027
$ws.Cells.Item(19, 9).Value = 0.0006073039498957697
$ws.Cells.Item(19, 10).Value = 0.0006073039498957697
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 2.423449666666667
$ws.Cells.Item(19, 14).Value = 7.270349
$ws.Cells.Item(19, 15).Value = 0.03003061794650479
$ws.Cells.Item(19, 16).Value = 0.03003061794650479
$ws.Cells.Item(19, 17).Value = 0.1187708447136667
$ws.Cells.Item(19, 18).Value = 1.068937602423
$ws.Cells.Item(19, 19).Value = 0.00001823771289672315
$ws.Cells.Item(19, 20).Value = 0.00001823771289672315

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 0.3066033333333333
$ws.Cells.Item(20, 8).Value = 0.91981
$ws.Cells.Item(20, 9).Value = 0.003799331049083692
$ws.Cells.Item(20, 10).Value = 0.003799331049083692
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 66.85695766666667
$ws.Cells.Item(20, 14).Value = 200.570873
$ws.Cells.Item(20, 15).Value = 0.8284701681115905
$ws.Cells.Item(20, 16).Value = 0.8284701681115904
$ws.Cells.Item(20, 17).Value = 20.49856607712556
$ws.Cells.Item(20, 18).Value = 184.48709469413
$ws.Cells.Item(20, 19).Value = 0.003147632432945952
$ws.Cells.Item(20, 20).Value = 0.003147632432945951

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 0.3066033333333333
$ws.Cells.Item(21, 8).Value = 0.91981
$ws.Cells.Item(21, 9).Value = 0.003799331049083692
$ws.Cells.Item(21, 10).Value = 0.003799331049083692
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 10.944004
$ws.Cells.Item(21, 14).Value = 32.83201200000001
$ws.Cells.Item(21, 15).Value = 0.1356146188837786
$ws.Cells.Item(21, 16).Value = 0.1356146188837786
$ws.Cells.Item(21, 17).Value = 3.355468106413334
$ws.Cells.Item(21, 18).Value = 30.19921295772001
$ws.Cells.Item(21, 19).Value = 0.0005152448322347917
$ws.Cells.Item(21, 20).Value = 0.0005152448322347915

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 0.3066033333333333
$ws.Cells.Item(22, 8).Value = 0.91981
$ws.Cells.Item(22, 9).Value = 0.003799331049083692
$ws.Cells.Item(22, 10).Value = 0.003799331049083692
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 0.049009
$ws.Cells.Item(22, 14).Value = 0.147027
$ws.Cells.Item(22, 15).Value = 0.0006073039498957697
$ws.Cells.Item(22, 16).Value = 0.0006073039498957697
$ws.Cells.Item(22, 17).Value = 0.01502632276333333
$ws.Cells.Item(22, 18).Value = 0.13523690487
$ws.Cells.Item(22, 19).Value = 0.000002307348753070165
$ws.Cells.Item(22, 20).Value = 0.000002307348753070164

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 0.3066033333333333
$ws.Cells.Item(23, 8).Value = 0.91981
$ws.Cells.Item(23, 9).Value = 0.003799331049083692
$ws.Cells.Item(23, 10).Value = 0.003799331049083692
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 0.3066033333333333
$ws.Cells.Item(23, 14).Value = 0.91981
$ws.Cells.Item(23, 15).Value = 0.003799331049083692
$ws.Cells.Item(23, 16).Value = 0.003799331049083692
$ws.Cells.Item(23, 17).Value = 0.09400560401111112
$ws.Cells.Item(23, 18).Value = 0.8460504361000001
$ws.Cells.Item(23, 19).Value = 0.00001443491642053139
$ws.Cells.Item(23, 20).Value = 0.00001443491642053138

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 0.3066033333333333
$ws.Cells.Item(24, 8).Value = 0.91981
$ws.Cells.Item(24, 9).Value = 0.003799331049083692
$ws.Cells.Item(24, 10).Value = 0.003799331049083692
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 0.1192703333333334
$ws.Cells.Item(24, 14).Value = 0.357811
$ws.Cells.Item(24, 15).Value = 0.001477960059146655
$ws.Cells.Item(24, 16).Value = 0.001477960059146655
$ws.Cells.Item(24, 17).Value = 0.03656868176777778
$ws.Cells.Item(24, 18).Value = 0.3291181359100001
$ws.Cells.Item(24, 19).Value = 0.000005615259542021458
$ws.Cells.Item(24, 20).Value = 0.000005615259542021456

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 0.3066033333333333
$ws.Cells.Item(25, 8).Value = 0.91981
$ws.Cells.Item(25, 9).Value = 0.003799331049083692
$ws.Cells.Item(25, 10).Value = 0.003799331049083692
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 2.423449666666667
$ws.Cells.Item(25, 14).Value = 7.270349
$ws.Cells.Item(25, 15).Value = 0.03003061794650479
$ws.Cells.Item(25, 16).Value = 0.03003061794650479
$ws.Cells.Item(25, 17).Value = 0.7430377459655555
$ws.Cells.Item(25, 18).Value = 6.687339713689999
$ws.Cells.Item(25, 19).Value = 0.0001140962591873256
$ws.Cells.Item(25, 20).Value = 0.0001140962591873256

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 0.1192703333333334
$ws.Cells.Item(26, 8).Value = 0.357811
$ws.Cells.Item(26, 9).Value = 0.001477960059146655
$ws.Cells.Item(26, 10).Value = 0.001477960059146655
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 13).Value = 66.85695766666667
$ws.Cells.Item(26, 14).Value = 200.570873
$ws.Cells.Item(26, 15).Value = 0.8284701681115905
$ws.Cells.Item(26, 16).Value = 0.8284701681115904
$ws.Cells.Item(26, 17).Value = 7.974051626555891
$ws.Cells.Item(26, 18).Value = 71.766464639003
$ws.Cells.Item(26, 19).Value = 0.001224445818663446
$ws.Cells.Item(26, 20).Value = 0.001224445818663445

$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 0.1192703333333334
$ws.Cells.Item(27, 8).Value = 0.357811
$ws.Cells.Item(27, 9).Value = 0.001477960059146655
$ws.Cells.Item(27, 10).Value = 0.001477960059146655
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 10.944004
$ws.Cells.Item(27, 14).Value = 32.83201200000001
$ws.Cells.Item(27, 15).Value = 0.1356146188837786
$ws.Cells.Item(27, 16).Value = 0.1356146188837786
$ws.Cells.Item(27, 17).Value = 1.305295005081334
$ws.Cells.Item(27, 18).Value = 11.747655045732
$ws.Cells.Item(27, 19).Value = 0.0002004329901466205
$ws.Cells.Item(27, 20).Value = 0.0002004329901466205

$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 0.1192703333333334
$ws.Cells.Item(28, 8).Value = 0.357811
$ws.Cells.Item(28, 9).Value = 0.001477960059146655
$ws.Cells.Item(28, 10).Value = 0.001477960059146655
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 0.049009
$ws.Cells.Item(28, 14).Value = 0.147027
$ws.Cells.Item(28, 15).Value = 0.0006073039498957697
$ws.Cells.Item(28, 16).Value = 0.0006073039498957697
$ws.Cells.Item(28, 17).Value = 0.005845319766333334
$ws.Cells.Item(28, 18).Value = 0.05260787789700001
$ws.Cells.Item(28, 19).Value = 0.0000008975709817079493
$ws.Cells.Item(28, 20).Value = 0.000000897570981707949

$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 0.1192703333333334
$ws.Cells.Item(29, 8).Value = 0.357811
$ws.Cells.Item(29, 9).Value = 0.001477960059146655
$ws.Cells.Item(29, 10).Value = 0.001477960059146655
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 0.3066033333333333
$ws.Cells.Item(29, 14).Value = 0.91981
$ws.Cells.Item(29, 15).Value = 0.003799331049083692
$ws.Cells.Item(29, 16).Value = 0.003799331049083692
$ws.Cells.Item(29, 17).Value = 0.03656868176777778
$ws.Cells.Item(29, 18).Value = 0.3291181359100001
$ws.Cells.Item(29, 19).Value = 0.000005615259542021458
$ws.Cells.Item(29, 20).Value = 0.000005615259542021456

$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 0.1192703333333334
$ws.Cells.Item(30, 8).Value = 0.357811
$ws.Cells.Item(30, 9).Value = 0.001477960059146655
$ws.Cells.Item(30, 10).Value = 0.001477960059146655
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 0.1192703333333334
$ws.Cells.Item(30, 14).Value = 0.357811
$ws.Cells.Item(30, 15).Value = 0.001477960059146655
$ws.Cells.Item(30, 16).Value = 0.001477960059146655
$ws.Cells.Item(30, 17).Value = 0.01422541241344445
$ws.Cells.Item(30, 18).Value = 0.128028711721
$ws.Cells.Item(30, 19).Value = 0.000002184365936432785
$ws.Cells.Item(30, 20).Value = 0.000002184365936432785

$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 0.1192703333333334
$ws.Cells.Item(31, 8).Value = 0.357811
$ws.Cells.Item(31, 9).Value = 0.001477960059146655
$ws.Cells.Item(31, 10).Value = 0.001477960059146655
$ws.Cells.Item(31, 11).Value = 3
$ws.Cells.Item(31, 13).Value = 2.423449666666667
$ws.Cells.Item(31, 14).Value = 7.270349
$ws.Cells.Item(31, 15).Value = 0.03003061794650479
$ws.Cells.Item(31, 16).Value = 0.03003061794650479
$ws.Cells.Item(31, 17).Value = 0.2890456495598889
$ws.Cells.Item(31, 18).Value = 2.601410846039
$ws.Cells.Item(31, 19).Value = 0.00004438405387642683
$ws.Cells.Item(31, 20).Value = 0.00004438405387642682

$ws.Cells.Item(32, 5).Value = 3
$ws.Cells.Item(32, 7).Value = 2.423449666666667
$ws.Cells.Item(32, 8).Value = 7.270349
$ws.Cells.Item(32, 9).Value = 0.03003061794650479
$ws.Cells.Item(32, 10).Value = 0.03003061794650479
$ws.Cells.Item(32, 11).Value = 3
$ws.Cells.Item(32, 13).Value = 66.85695766666667
$ws.Cells.Item(32, 14).Value = 200.570873
$ws.Cells.Item(32, 15).Value = 0.8284701681115905
$ws.Cells.Item(32, 16).Value = 0.8284701681115904
$ws.Cells.Item(32, 17).Value = 162.0244717716308
$ws.Cells.Item(32, 18).Value = 1458.220245944677
$ws.Cells.Item(32, 19).Value = 0.02487947109863577
$ws.Cells.Item(32, 20).Value = 0.02487947109863576

$ws.Cells.Item(33, 5).Value = 3
$ws.Cells.Item(33, 7).Value = 2.423449666666667
$ws.Cells.Item(33, 8).Value = 7.270349
$ws.Cells.Item(33, 9).Value = 0.03003061794650479
$ws.Cells.Item(33, 10).Value = 0.03003061794650479
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 10.944004
$ws.Cells.Item(33, 14).Value = 32.83201200000001
$ws.Cells.Item(33, 15).Value = 0.1356146188837786
$ws.Cells.Item(33, 16).Value = 0.1356146188837786
$ws.Cells.Item(33, 17).Value = 26.52224284579867
$ws.Cells.Item(33, 18).Value = 238.700185612188
$ws.Cells.Item(33, 19).Value = 0.004072590807659609
$ws.Cells.Item(33, 20).Value = 0.004072590807659608

$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 7).Value = 2.423449666666667
$ws.Cells.Item(34, 8).Value = 7.270349
$ws.Cells.Item(34, 9).Value = 0.03003061794650479
$ws.Cells.Item(34, 10).Value = 0.03003061794650479
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 0.049009
$ws.Cells.Item(34, 14).Value = 0.147027
$ws.Cells.Item(34, 15).Value = 0.0006073039498957697
$ws.Cells.Item(34, 16).Value = 0.0006073039498957697
$ws.Cells.Item(34, 17).Value = 0.1187708447136667
$ws.Cells.Item(34, 18).Value = 1.068937602423
$ws.Cells.Item(34, 19).Value = 0.00001823771289672315
$ws.Cells.Item(34, 20).Value = 0.00001823771289672315

$ws.Cells.Item(35, 5).Value = 3
$ws.Cells.Item(35, 7).Value = 2.423449666666667
$ws.Cells.Item(35, 8).Value = 7.270349
$ws.Cells.Item(35, 9).Value = 0.03003061794650479
$ws.Cells.Item(35, 10).Value = 0.03003061794650479
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 0.3066033333333333
$ws.Cells.Item(35, 14).Value = 0.91981
$ws.Cells.Item(35, 15).Value = 0.003799331049083692
$ws.Cells.Item(35, 16).Value = 0.003799331049083692
$ws.Cells.Item(35, 17).Value = 0.7430377459655555
$ws.Cells.Item(35, 18).Value = 6.687339713689999
$ws.Cells.Item(35, 19).Value = 0.0001140962591873256
$ws.Cells.Item(35, 20).Value = 0.0001140962591873256

$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 7).Value = 2.423449666666667
$ws.Cells.Item(36, 8).Value = 7.270349
$ws.Cells.Item(36, 9).Value = 0.03003061794650479
$ws.Cells.Item(36, 10).Value = 0.03003061794650479
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 0.1192703333333334
$ws.Cells.Item(36, 14).Value = 0.357811
$ws.Cells.Item(36, 15).Value = 0.001477960059146655
$ws.Cells.Item(36, 16).Value = 0.001477960059146655
$ws.Cells.Item(36, 17).Value = 0.2890456495598889
$ws.Cells.Item(36, 18).Value = 2.601410846039
$ws.Cells.Item(36, 19).Value = 0.00004438405387642683
$ws.Cells.Item(36, 20).Value = 0.00004438405387642682

$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 7).Value = 2.423449666666667
$ws.Cells.Item(37, 8).Value = 7.270349
$ws.Cells.Item(37, 9).Value = 0.03003061794650479
$ws.Cells.Item(37, 10).Value = 0.03003061794650479
$ws.Cells.Item(37, 11).Value = 3
$ws.Cells.Item(37, 13).Value = 2.423449666666667
$ws.Cells.Item(37, 14).Value = 7.270349
$ws.Cells.Item(37, 15).Value = 0.03003061794650479
$ws.Cells.Item(37, 16).Value = 0.03003061794650479
$ws.Cells.Item(37, 17).Value = 5.873108286866777
$ws.Cells.Item(37, 18).Value = 52.857974581801
$ws.Cells.Item(37, 19).Value = 0.0009018380142489355
$ws.Cells.Item(37, 20).Value = 0.0009018380142489353
